$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 14) describing "Muse 125 ug, Urethrastab" package,
# inserted right after the current last data row (row 13).

# Copy row 13's formatting (styles/number formats) down into row 14 first,
# so the new row visually matches the rest of the table.
$ws.Range("A13:R13").Copy()
$ws.Range("A14:R14").PasteSpecial(-4122)
$ws.Rows.Item(14).RowHeight = 12.75

# Now fill in the actual values for the new row.
$ws.Cells.Item(14, 1).Value2 = 54525                  # Zulassungs-Nummer
$ws.Cells.Item(14, 2).Value2 = 1                      # Dosistärke-nummer
$ws.Cells.Item(14, 3).Value2 = "Muse 125 ug, Urethrastab"          # Präparatebezeichnung
$ws.Cells.Item(14, 4).Value2 = "MEDA Pharma GmbH"                  # Zulassungsinhaberin
$ws.Cells.Item(14, 5).Value2 = "05.99.0."                          # IT-Nummer
$ws.Cells.Item(14, 6).Value2 = "G04BE01"                           # ATC-Code
$ws.Cells.Item(14, 7).Value2 = $ws.Cells.Item(13, 7).Value2        # Heilmittelcode (reuse "Synthetika human")
$ws.Cells.Item(14, 8).Value2 = 35852                  # Erstzul.datum Präp.
$ws.Cells.Item(14, 9).Value2 = 35852                  # Zul.datum Dosisstärke *
$ws.Cells.Item(14, 10).Value2 = 43553                 # Gültigkeits-datum *
$ws.Cells.Item(14, 11).Value2 = 36                    # Verpackungs ID
$ws.Cells.Item(14, 12).Value2 = "6x1 Urethrastab"                  # Packungsgrösse
# Column 13 (Einheit) stays blank, same as row 13.
$ws.Cells.Item(14, 14).Value2 = $ws.Cells.Item(13, 14).Value2      # Abgabekategorie (reuse "B")
$ws.Cells.Item(14, 15).Value2 = "alprostadilum"                    # Wirkstoff
$ws.Cells.Item(14, 16).Value2 = "alprostadilum 125 µg, excipiens ad gelatum pro praeparatione."  # Zusammensetzung
$ws.Cells.Item(14, 17).Value2 = "Erektile Dysfunktion"             # Anwendungsgebiet Präparate
# Column 18 (Anwendungsgebiet Dosisstärke) stays blank, same as row 13.

# Move the view/selection down onto the newly added row, matching Excel's
# behaviour of following the last edited row.
$ws.Range("A14").Select()

"done"
